# OLX Monitor run 2026-02-22 13:31 — append the latest monitoring snapshot
# for the "wszystkie-lublin" profile onto the PODSUMOWANIE log sheet.
#
# The new check re-observed the same 8 listings as the previous run
# (2026-02-22 12:40:56, rows 139-146) so the new rows (147-154) are
# identical in every column except the "checked at" timestamp in col A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

$lastRow = 146
$numRows = 8
$newFirstRow = $lastRow + 1
$newLastRow = $lastRow + $numRows

$srcRange = "A" + ($lastRow - $numRows + 1) + ":H" + $lastRow
$dstRange = "A" + $newFirstRow + ":H" + $newLastRow

# Duplicate the previous batch of 8 listing rows (values + number formats +
# styles) into the new block of rows right below the existing data.
$ws.Range($srcRange).Copy()
$ws.Range($dstRange).PasteSpecial(-4122)
$ws.Range($srcRange).Copy()
$ws.Range($dstRange).PasteSpecial(-4104)

# Stamp the new rows with this run's "checked at" timestamp.
$ws.Range("A" + $newFirstRow + ":A" + $newLastRow).Value = "2026-02-22 13:31:28"
